$d = $word.ActiveDocument

# Find the paragraph that contains "Missing variable :" and remove the
# trailing error-marker runs ("    <---M2Doc version mismatch: ... 3.2.0    ")
# leaving only "Missing variable :" in the paragraph.

$para = $d.Paragraphs(1).Range

$rng = $para.Duplicate
$found = $rng.Find.Execute(
    "Missing variable :",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    # Move to the end of the found text, then extend to the end of the
    # paragraph (excluding the paragraph mark) and delete that range.
    $delStart = $rng.End
    $paraEnd = $para.End - 1  # exclude paragraph mark
    if ($paraEnd -gt $delStart) {
        $delRange = $d.Range($delStart, $paraEnd)
        $delRange.Delete()
    }
}
